$wb = $excel.ActiveWorkbook

# Sheet 1: ALC
$ws = $wb.Worksheets.Item(1)
$ws.Range("H31").Value = 163.2
$ws.Range("I31").Value = 163.2
$ws.Range("K31").Value = 489.6
$ws.Range("M31").Value = -259.6
$ws.Range("H43").Value = 885.3125
$ws.Range("I43").Value = 724.875
$ws.Range("J43").Value = 1045.75
$ws.Range("K43").Value = 724.875
$ws.Range("L43").Value = 1045.75
$ws.Range("M43").Value = -655.875
$ws.Range("N43").Value = -1183.75
$ws.Range("H51").Value = 4280.1177
$ws.Range("I51").Value = 1671.4286
$ws.Range("J51").Value = 6106.2
$ws.Range("K51").Value = 1671.4286
$ws.Range("L51").Value = 6106.2
$ws.Range("M51").Value = -1187.4286
$ws.Range("N51").Value = -7074.2
$ws.Range("H53").Value = 1236.625
$ws.Range("I53").Value = 1735.8182
$ws.Range("K53").Value = 1735.8182
$ws.Range("M53").Value = -1098.8182
$ws.Range("H74").Value = 3211.5
$ws.Range("I74").Value = 2718
$ws.Range("J74").Value = 3435.818
$ws.Range("K74").Value = 2718
$ws.Range("L74").Value = 3435.818
$ws.Range("M74").Value = -1782
$ws.Range("N74").Value = -5307.818
$ws.Range("H76").Value = 2950
$ws.Range("I76").Value = 2950
$ws.Range("K76").Value = 2950
$ws.Range("M76").Value = -2635
$ws.Range("H77").Value = 3211.5
$ws.Range("I77").Value = 2718
$ws.Range("J77").Value = 3435.818
$ws.Range("K77").Value = 13590
$ws.Range("L77").Value = 17179.09
$ws.Range("M77").Value = -8910
$ws.Range("N77").Value = -26539.09
$ws.Range("H79").Value = 2950
$ws.Range("I79").Value = 2950
$ws.Range("K79").Value = 2950
$ws.Range("M79").Value = -1858
$ws.Range("H88").Value = 3216.05
$ws.Range("I88").Value = 3126.125
$ws.Range("J88").Value = 3276
$ws.Range("K88").Value = 3126.125
$ws.Range("L88").Value = 3276
$ws.Range("M88").Value = -2720.125
$ws.Range("N88").Value = -4088
$ws.Range("H91").Value = 3216.05
$ws.Range("I91").Value = 3126.125
$ws.Range("J91").Value = 3276
$ws.Range("K91").Value = 3126.125
$ws.Range("L91").Value = 3276
$ws.Range("M91").Value = -1722.125
$ws.Range("N91").Value = -6084
$ws.Range("H113").Value = 4673.4736
$ws.Range("I113").Value = 3520
$ws.Range("J113").Value = 5955.1113
$ws.Range("K113").Value = 3520
$ws.Range("L113").Value = 5955.1113
$ws.Range("M113").Value = -266
$ws.Range("N113").Value = -12463.1113
$ws.Range("H138").Value = 3776734
$ws.Range("I138").Value = 2292.2144
$ws.Range("J138").Value = 5131662
$ws.Range("K138").Value = 6876.6432
$ws.Range("L138").Value = 15394986
$ws.Range("M138").Value = -1736.6432
$ws.Range("N138").Value = -15405266
$ws.Range("H139").Value = 39898.25
$ws.Range("J139").Value = 39898.25
$ws.Range("L139").Value = 39898.25
$ws.Range("N139").Value = -50178.25
$ws.Range("H141").Value = 1024.6923
$ws.Range("I141").Value = 691.8929000000001
$ws.Range("J141").Value = 1871.8182
$ws.Range("K141").Value = 2075.6787
$ws.Range("L141").Value = 5615.4546
$ws.Range("M141").Value = 3104.3213
$ws.Range("N141").Value = -15975.4546

# Sheet 2: ARM
$ws = $wb.Worksheets.Item(2)
$ws.Range("H74").Value = 23872.25
$ws.Range("I74").Value = 36686.57
$ws.Range("K74").Value = 36686.57
$ws.Range("M74").Value = -35812.57
$ws.Range("H77").Value = 23872.25
$ws.Range("I77").Value = 36686.57
$ws.Range("K77").Value = 183432.85
$ws.Range("M77").Value = -179064.85
$ws.Range("H102").Value = 1457.3704
$ws.Range("I102").Value = 1252.1666
$ws.Range("K102").Value = 1252.1666
$ws.Range("M102").Value = 369.8334

# Sheet 3: BSM
$ws = $wb.Worksheets.Item(3)
$ws.Range("H82").Value = 9113.75
$ws.Range("I82").Value = 6130
$ws.Range("J82").Value = 30000
$ws.Range("K82").Value = 6130
$ws.Range("L82").Value = 30000
$ws.Range("M82").Value = -5747
$ws.Range("N82").Value = -30766
$ws.Range("H85").Value = 9113.75
$ws.Range("I85").Value = 6130
$ws.Range("J85").Value = 30000
$ws.Range("K85").Value = 6130
$ws.Range("L85").Value = 30000
$ws.Range("M85").Value = -4804
$ws.Range("N85").Value = -32652
$ws.Range("H86").Value = 2318
$ws.Range("I86").Value = 2244
$ws.Range("K86").Value = 2244
$ws.Range("M86").Value = -1121
$ws.Range("H89").Value = 2318
$ws.Range("I89").Value = 2244
$ws.Range("K89").Value = 11220
$ws.Range("M89").Value = -5604
$ws.Range("H99").Value = 1649.5
$ws.Range("I99").Value = 1561.3334
$ws.Range("K99").Value = 1561.3334
$ws.Range("M99").Value = -63.33339999999998
$ws.Range("H134").Value = 600445.75
$ws.Range("I134").Value = 932750.9
$ws.Range("J134").Value = 5065.7915
$ws.Range("K134").Value = 2798252.7
$ws.Range("L134").Value = 15197.3745
$ws.Range("M134").Value = -2795717.7
$ws.Range("N134").Value = -20267.3745

# Sheet 4: CRP
$ws = $wb.Worksheets.Item(4)
$ws.Range("H62").Value = 2299.3125
$ws.Range("I62").Value = 2262.4167
$ws.Range("J62").Value = 2410
$ws.Range("K62").Value = 2262.4167
$ws.Range("L62").Value = 2410
$ws.Range("M62").Value = -1638.4167
$ws.Range("N62").Value = -3658
$ws.Range("H65").Value = 2299.3125
$ws.Range("I65").Value = 2262.4167
$ws.Range("J65").Value = 2410
$ws.Range("K65").Value = 11312.0835
$ws.Range("L65").Value = 12050
$ws.Range("M65").Value = -8192.083500000001
$ws.Range("N65").Value = -18290
$ws.Range("H134").Value = 18645676
$ws.Range("I134").Value = 2001524.9
$ws.Range("J134").Value = 111113176
$ws.Range("K134").Value = 6004574.699999999
$ws.Range("L134").Value = 333339528
$ws.Range("M134").Value = -6002039.699999999
$ws.Range("N134").Value = -333344598

# Sheet 6: GSM
$ws = $wb.Worksheets.Item(6)
$ws.Range("H80").Value = 2459.2856
$ws.Range("I80").Value = 2350.3572
$ws.Range("J80").Value = 2677.1428
$ws.Range("K80").Value = 2350.3572
$ws.Range("L80").Value = 2677.1428
$ws.Range("M80").Value = -1352.3572
$ws.Range("N80").Value = -4673.1428
$ws.Range("H83").Value = 2459.2856
$ws.Range("I83").Value = 2350.3572
$ws.Range("J83").Value = 2677.1428
$ws.Range("K83").Value = 11751.786
$ws.Range("L83").Value = 13385.714
$ws.Range("M83").Value = -6759.786
$ws.Range("N83").Value = -23369.714
$ws.Range("H132").Value = 2383414.2
$ws.Range("I132").Value = 2501.64
$ws.Range("J132").Value = 5884756
$ws.Range("K132").Value = 7504.92
$ws.Range("L132").Value = 17654268
$ws.Range("M132").Value = -4974.92
$ws.Range("N132").Value = -17659328

# Sheet 7: LTW
$ws = $wb.Worksheets.Item(7)
$ws.Range("H82").Value = 1194.6
$ws.Range("J82").Value = 950
$ws.Range("L82").Value = 950
$ws.Range("N82").Value = -1672
$ws.Range("H85").Value = 1194.6
$ws.Range("J85").Value = 950
$ws.Range("L85").Value = 950
$ws.Range("N85").Value = -3446
$ws.Range("H116").Value = 32000
$ws.Range("J116").Value = 32000
$ws.Range("L116").Value = 32000
$ws.Range("N116").Value = -41178
$ws.Range("H138").Value = 33419.25
$ws.Range("J138").Value = 33419.25
$ws.Range("L138").Value = 33419.25
$ws.Range("N138").Value = -43699.25
$ws.Range("H139").Value = 0
$ws.Range("J139").Value = 0
$ws.Range("L139").ClearContents() | Out-Null
$ws.Range("N139").ClearContents() | Out-Null

# Sheet 8: WVR
$ws = $wb.Worksheets.Item(8)
$ws.Range("H27").Value = 30000
$ws.Range("J27").Value = 30000
$ws.Range("L27").Value = 30000
$ws.Range("N27").Value = -30138
$ws.Range("H115").Value = 31344.25
$ws.Range("J115").Value = 31344.25
$ws.Range("L115").Value = 31344.25
$ws.Range("N115").Value = -34478.25
$ws.Range("H132").Value = 2019.6818
$ws.Range("I132").Value = 2002.3726
$ws.Range("J132").Value = 2078.5334
$ws.Range("K132").Value = 6007.1178
$ws.Range("L132").Value = 6235.600199999999
$ws.Range("M132").Value = -3477.1178
$ws.Range("N132").Value = -11295.6002
